# Applies the "Add files via upload" edit to the Partidos sheet:
#  1. Flips 50 existing cells in columns I/J (tarjetas_amarillas /
#     tarjetas_rojas) from 1 back to 0 - these were apparently a paste
#     error that got corrected upstream.
#  2. Appends 10 new match rows (351-360) for the 2025-09-20 (serial
#     45864) session.
#  3. Updates the sheet selection to reflect where the author ended up
#     after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partidos")
$ws.Activate()

# ---------------------------------------------------------------------
# 1. Cells that changed from 1 -> 0
# ---------------------------------------------------------------------
$iRows = @(7,36,37,41,43,46,47,52,59,68,81,82,83,93,95,100,101,105,110,111,112,116,117,118,119,120,121,147,156,161,171,181,186,216,222,227,229,236,248,267,276,282,290,291,296,301)
foreach ($r in $iRows) {
    $ws.Cells.Item($r, 9).Value = 0   # column I = tarjetas_amarillas
}

$jRows = @(122,237,238,311)
foreach ($r in $jRows) {
    $ws.Cells.Item($r, 10).Value = 0  # column J = tarjetas_rojas
}

# ---------------------------------------------------------------------
# 2. New rows 351-360
#    columns: A fecha, B jugador, C equipo, D posicion, E goles,
#             F autogoles, G arquero, H goles_recibidos,
#             I tarjetas_amarillas, J tarjetas_rojas, K asistencias,
#             L Penales_Atajados
# ---------------------------------------------------------------------
$newRows = @(
    @(45864, "Gember Marin Sarria",        "Azul",     "Arquero",       0, 0, $true,  2, 0, 0, 0, 0),
    @(45864, "Edwin Casas",                "Amarillo", "Arquero",       0, 0, $true,  2, 0, 0, 0, 0),
    @(45864, "Arnul David Narvaez",        "Azul",     "Delantero",     1, 0, $false, 0, 0, 1, 0, 0),
    @(45864, "Jefferson Delgado",          "Azul",     "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0),
    @(45864, "Cesar Augusto Estrada",      "Azul",     "Delantero",     1, 0, $false, 0, 0, 0, 0, 0),
    @(45864, "Andres Tangarife",           "Azul",     "Delantero",     0, 0, $false, 0, 0, 0, 1, 0),
    @(45864, "Sebastian Giraldo",          "Amarillo", "Mediocampista", 1, 0, $false, 0, 0, 0, 1, 0),
    @(45864, "Andres Jurado",              "Amarillo", "Delantero",     1, 0, $false, 0, 0, 0, 0, 0),
    @(45864, "Jose Antonio Nieva Chaves",  "Amarillo", "Defensa",       0, 0, $false, 0, 1, 0, 0, 0),
    @(45864, "Juan Felipe Gutierrez",      "Amarillo", "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0)
)

$startRow = 351
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 5).Value = $data[4]
    $ws.Cells.Item($row, 6).Value = $data[5]
    $ws.Cells.Item($row, 7).Value = $data[6]
    $ws.Cells.Item($row, 8).Value = $data[7]
    $ws.Cells.Item($row, 9).Value = $data[8]
    $ws.Cells.Item($row, 10).Value = $data[9]
    $ws.Cells.Item($row, 11).Value = $data[10]
    $ws.Cells.Item($row, 12).Value = $data[11]
}

# ---------------------------------------------------------------------
# 3. Final selection state (matches the author's ending cursor position)
# ---------------------------------------------------------------------
$ws.Range("B362").Select()
